$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8914.0
$ws.Range("D2").Value = 50.647727966308594

$ws.Range("B4").Value = 11715.0
$ws.Range("D4").Value = 33.28125

$ws.Range("B6").Value = 78.0
$ws.Range("D6").Value = 97.5
